# UniDebAutomCar2017-1_Requirements.xlsx - Team1 update
# - Added a note in column H (Comment) for the UNIDEB_4 row explaining
#   that the following 6 requirements are needed from Team1 for Sprint1.
# - The AutoFilter range was shrunk back down to just the header row.
# - The underlying _FilterDatabase defined names were refreshed to track
#   the new filter range (mirroring what the spreadsheet application does
#   whenever the AutoFilter range changes).
# - Selection moved onto the newly edited cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirement")

# 1. Add the new comment text for UNIDEB_4 row (row 6), column H ("Comment")
$ws.Range("H6").Value = "The following 6 requirements are needed from Team1 for Sprint1"

# 2. Shrink the AutoFilter back to just the header row (A2:H2)
$ws.AutoFilterMode = $false
$ws.Range("A2:H2").AutoFilter()

# 3. Refresh the _FilterDatabase defined names to reflect the new filter range
$names = $wb.Names
$names.Item(1).RefersTo = "=Requirement!`$A`$2:`$H`$2"
$names.Item(2).RefersTo = "=Requirement!`$A`$2:`$H`$12"
$names.Item(3).RefersTo = "=Requirement!`$A`$2:`$H`$2"
$names.Item(4).RefersTo = "=Requirement!`$A`$2:`$H`$12"
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0", "=Requirement!`$A`$2:`$H`$2")

# 4. Reflect the updated selection/scroll position used when the edit was made
$ws.Range("B1").Select()
$ws.Range("H6").Select()
